$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (ALC)
$ws.Cells.Item(19, 8).Value = 4037.8333  # H19: 4087.8333 -> 4037.8333
$ws.Cells.Item(19, 9).Value = 10375  # I19: 7250 -> 10375
$ws.Cells.Item(19, 10).Value = 869.25  # J19: 925.6667 -> 869.25
$ws.Cells.Item(19, 11).Value = 10375  # K19: 7250 -> 10375
$ws.Cells.Item(19, 12).Value = 869.25  # L19: 925.6667 -> 869.25
$ws.Cells.Item(19, 13).Value = -10200  # M19: -7075 -> -10200
$ws.Cells.Item(19, 14).Value = -1219.25  # N19: -1275.6667 -> -1219.25

# Row 43 (ALC)
$ws.Cells.Item(43, 8).Value = 666.6667  # H43: 900 -> 666.6667
$ws.Cells.Item(43, 9).Value = 750  # I43: 900 -> 750
$ws.Cells.Item(43, 10).Value = 500  # J43: 0 -> 500
$ws.Cells.Item(43, 11).Value = 750  # K43: 900 -> 750
$ws.Cells.Item(43, 12).Value = 500  # L43: 0 -> 500
$ws.Cells.Item(43, 13).Value = -681  # M43: -831 -> -681
$ws.Cells.Item(43, 14).Value = -638  # N43: None -> -638

# Row 53 (ALC)
$ws.Cells.Item(53, 8).Value = 3829.9092  # H53: 3792.6365 -> 3829.9092
$ws.Cells.Item(53, 9).Value = 245.5  # I53: 143 -> 245.5
$ws.Cells.Item(53, 11).Value = 245.5  # K53: 143 -> 245.5
$ws.Cells.Item(53, 13).Value = 391.5  # M53: 494 -> 391.5

# Row 116 (ALC)
$ws.Cells.Item(116, 8).Value = 5299.615  # H116: 5338.077 -> 5299.615
$ws.Cells.Item(116, 9).Value = 2899.8  # I116: 2999.8 -> 2899.8
$ws.Cells.Item(116, 11).Value = 2899.8  # K116: 2999.8 -> 2899.8
$ws.Cells.Item(116, 13).Value = 542.1999999999998  # M116: 442.1999999999998 -> 542.1999999999998

# Row 129 (ALC)
$ws.Cells.Item(129, 8).Value = 179522.72  # H129: 182793.69 -> 179522.72
$ws.Cells.Item(129, 9).Value = 281.83334  # I129: 324.25 -> 281.83334
$ws.Cells.Item(129, 10).Value = 201031.62  # J129: 197105.02 -> 201031.62
$ws.Cells.Item(129, 11).Value = 845.5000200000001  # K129: 972.75 -> 845.5000200000001
$ws.Cells.Item(129, 12).Value = 603094.86  # L129: 591315.0599999999 -> 603094.86
$ws.Cells.Item(129, 13).Value = 4154.49998  # M129: 4027.25 -> 4154.49998
$ws.Cells.Item(129, 14).Value = -613094.86  # N129: -601315.0599999999 -> -613094.86

# Row 132 (ALC)
$ws.Cells.Item(132, 8).Value = 2572.919  # H132: 2573.8647 -> 2572.919
$ws.Cells.Item(132, 9).Value = 2602.5557  # I132: 2603.5278 -> 2602.5557
$ws.Cells.Item(132, 11).Value = 7807.6671  # K132: 7810.5834 -> 7807.6671
$ws.Cells.Item(132, 13).Value = -5277.6671  # M132: -5280.5834 -> -5277.6671

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 1379.7142  # H137: 1305.2258 -> 1379.7142
$ws.Cells.Item(137, 9).Value = 1397.0435  # I137: 1281.75 -> 1397.0435
$ws.Cells.Item(137, 10).Value = 1300  # J137: 1385.7142 -> 1300
$ws.Cells.Item(137, 11).Value = 4191.1305  # K137: 3845.25 -> 4191.1305
$ws.Cells.Item(137, 12).Value = 3900  # L137: 4157.142599999999 -> 3900
$ws.Cells.Item(137, 13).Value = -1641.1305  # M137: -1295.25 -> -1641.1305
$ws.Cells.Item(137, 14).Value = -9000  # N137: -9257.142599999999 -> -9000

# Row 141 (ALC)
$ws.Cells.Item(141, 8).Value = 2134  # H141: 2677.6 -> 2134
$ws.Cells.Item(141, 9).Value = 1771.8462  # I141: 2295.3635 -> 1771.8462
$ws.Cells.Item(141, 10).Value = 3703.3333  # J141: 3728.75 -> 3703.3333
$ws.Cells.Item(141, 11).Value = 5315.5386  # K141: 6886.0905 -> 5315.5386
$ws.Cells.Item(141, 12).Value = 11109.9999  # L141: 11186.25 -> 11109.9999
$ws.Cells.Item(141, 13).Value = -135.5385999999999  # M141: -1706.0905 -> -135.5385999999999
$ws.Cells.Item(141, 14).Value = -21469.9999  # N141: -21546.25 -> -21469.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Cells.Item(32, 8).Value = 2494.85  # H32: 2841.57 -> 2494.85
$ws.Cells.Item(32, 9).Value = 2226.4575  # I32: 2401.7283 -> 2226.4575
$ws.Cells.Item(32, 10).Value = 6699.6665  # J32: 7899.75 -> 6699.6665
$ws.Cells.Item(32, 11).Value = 2226.4575  # K32: 2401.7283 -> 2226.4575
$ws.Cells.Item(32, 12).Value = 6699.6665  # L32: 7899.75 -> 6699.6665
$ws.Cells.Item(32, 13).Value = -1939.4575  # M32: -2114.7283 -> -1939.4575
$ws.Cells.Item(32, 14).Value = -7273.6665  # N32: -8473.75 -> -7273.6665

# Row 61 (ARM)
$ws.Cells.Item(61, 8).Value = 2252.2173  # H61: 2295.0386 -> 2252.2173
$ws.Cells.Item(61, 9).Value = 1782.5625  # I61: 1778.4736 -> 1782.5625
$ws.Cells.Item(61, 10).Value = 3325.7144  # J61: 3697.1428 -> 3325.7144
$ws.Cells.Item(61, 11).Value = 1782.5625  # K61: 1778.4736 -> 1782.5625
$ws.Cells.Item(61, 12).Value = 3325.7144  # L61: 3697.1428 -> 3325.7144
$ws.Cells.Item(61, 13).Value = -1570.5625  # M61: -1566.4736 -> -1570.5625
$ws.Cells.Item(61, 14).Value = -3749.7144  # N61: -4121.1428 -> -3749.7144

# Row 74 (ARM)
$ws.Cells.Item(74, 8).Value = 45456668  # H74: 27779254 -> 45456668
$ws.Cells.Item(74, 9).Value = 83334010  # I74: 35714876 -> 83334010
$ws.Cells.Item(74, 10).Value = 3861.3  # J74: 4576.625 -> 3861.3
$ws.Cells.Item(74, 11).Value = 83334010  # K74: 35714876 -> 83334010
$ws.Cells.Item(74, 12).Value = 3861.3  # L74: 4576.625 -> 3861.3
$ws.Cells.Item(74, 13).Value = -83333136  # M74: -35714002 -> -83333136
$ws.Cells.Item(74, 14).Value = -5609.3  # N74: -6324.625 -> -5609.3

# Row 77 (ARM)
$ws.Cells.Item(77, 8).Value = 45456668  # H77: 27779254 -> 45456668
$ws.Cells.Item(77, 9).Value = 83334010  # I77: 35714876 -> 83334010
$ws.Cells.Item(77, 10).Value = 3861.3  # J77: 4576.625 -> 3861.3
$ws.Cells.Item(77, 11).Value = 416670050  # K77: 178574380 -> 416670050
$ws.Cells.Item(77, 12).Value = 19306.5  # L77: 22883.125 -> 19306.5
$ws.Cells.Item(77, 13).Value = -416665682  # M77: -178570012 -> -416665682
$ws.Cells.Item(77, 14).Value = -28042.5  # N77: -31619.125 -> -28042.5

# Row 102 (ARM)
$ws.Cells.Item(102, 8).Value = 1482.8667  # H102: 1606.5 -> 1482.8667
$ws.Cells.Item(102, 9).Value = 1380.1818  # I102: 1503.75 -> 1380.1818
$ws.Cells.Item(102, 10).Value = 1765.25  # J102: 1743.5 -> 1765.25
$ws.Cells.Item(102, 11).Value = 1380.1818  # K102: 1503.75 -> 1380.1818
$ws.Cells.Item(102, 12).Value = 1765.25  # L102: 1743.5 -> 1765.25
$ws.Cells.Item(102, 13).Value = 241.8181999999999  # M102: 118.25 -> 241.8181999999999
$ws.Cells.Item(102, 14).Value = -5009.25  # N102: -4987.5 -> -5009.25

# Row 136 (ARM)
$ws.Cells.Item(136, 8).Value = 2252.2173  # H136: 2295.0386 -> 2252.2173
$ws.Cells.Item(136, 9).Value = 1782.5625  # I136: 1778.4736 -> 1782.5625
$ws.Cells.Item(136, 10).Value = 3325.7144  # J136: 3697.1428 -> 3325.7144
$ws.Cells.Item(136, 11).Value = 5347.6875  # K136: 5335.4208 -> 5347.6875
$ws.Cells.Item(136, 12).Value = 9977.143199999999  # L136: 11091.4284 -> 9977.143199999999
$ws.Cells.Item(136, 13).Value = -2797.6875  # M136: -2785.4208 -> -2797.6875
$ws.Cells.Item(136, 14).Value = -15077.1432  # N136: -16191.4284 -> -15077.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 80 (BSM)
$ws.Cells.Item(80, 8).Value = 1200.7931  # H80: 904.12 -> 1200.7931
$ws.Cells.Item(80, 9).Value = 813.63635  # I80: 1053 -> 813.63635
$ws.Cells.Item(80, 10).Value = 1437.3889  # J80: 834.05884 -> 1437.3889
$ws.Cells.Item(80, 11).Value = 813.63635  # K80: 1053 -> 813.63635
$ws.Cells.Item(80, 12).Value = 1437.3889  # L80: 834.05884 -> 1437.3889
$ws.Cells.Item(80, 13).Value = 184.36365  # M80: -55 -> 184.36365
$ws.Cells.Item(80, 14).Value = -3433.3889  # N80: -2830.05884 -> -3433.3889

# Row 83 (BSM)
$ws.Cells.Item(83, 8).Value = 1200.7931  # H83: 904.12 -> 1200.7931
$ws.Cells.Item(83, 9).Value = 813.63635  # I83: 1053 -> 813.63635
$ws.Cells.Item(83, 10).Value = 1437.3889  # J83: 834.05884 -> 1437.3889
$ws.Cells.Item(83, 11).Value = 4068.18175  # K83: 5265 -> 4068.18175
$ws.Cells.Item(83, 12).Value = 7186.9445  # L83: 4170.2942 -> 7186.9445
$ws.Cells.Item(83, 13).Value = 923.8182500000003  # M83: -273 -> 923.8182500000003
$ws.Cells.Item(83, 14).Value = -17170.9445  # N83: -14154.2942 -> -17170.9445

# Row 134 (BSM)
$ws.Cells.Item(134, 8).Value = 3957.8572  # H134: 2941.45 -> 3957.8572
$ws.Cells.Item(134, 9).Value = 4104.88  # I134: 2920 -> 4104.88
$ws.Cells.Item(134, 10).Value = 2732.6667  # J134: 3349 -> 2732.6667
$ws.Cells.Item(134, 11).Value = 12314.64  # K134: 8760 -> 12314.64
$ws.Cells.Item(134, 12).Value = 8198.000100000001  # L134: 10047 -> 8198.000100000001
$ws.Cells.Item(134, 13).Value = -9779.639999999999  # M134: -6225 -> -9779.639999999999
$ws.Cells.Item(134, 14).Value = -13268.0001  # N134: -15117 -> -13268.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 3397.0977  # H31: 3747.7144 -> 3397.0977
$ws.Cells.Item(31, 9).Value = 2620.5217  # I31: 2981.1667 -> 2620.5217
$ws.Cells.Item(31, 10).Value = 4389.3887  # J31: 4559.353 -> 4389.3887
$ws.Cells.Item(31, 11).Value = 2620.5217  # K31: 2981.1667 -> 2620.5217
$ws.Cells.Item(31, 12).Value = 4389.3887  # L31: 4559.353 -> 4389.3887
$ws.Cells.Item(31, 13).Value = -2325.5217  # M31: -2686.1667 -> -2325.5217
$ws.Cells.Item(31, 14).Value = -4979.3887  # N31: -5149.353 -> -4979.3887

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 3397.0977  # H34: 3747.7144 -> 3397.0977
$ws.Cells.Item(34, 9).Value = 2620.5217  # I34: 2981.1667 -> 2620.5217
$ws.Cells.Item(34, 10).Value = 4389.3887  # J34: 4559.353 -> 4389.3887
$ws.Cells.Item(34, 11).Value = 2620.5217  # K34: 2981.1667 -> 2620.5217
$ws.Cells.Item(34, 12).Value = 4389.3887  # L34: 4559.353 -> 4389.3887
$ws.Cells.Item(34, 13).Value = -2418.5217  # M34: -2779.1667 -> -2418.5217
$ws.Cells.Item(34, 14).Value = -4793.3887  # N34: -4963.353 -> -4793.3887

# Row 62 (CRP)
$ws.Cells.Item(62, 8).Value = 52635740  # H62: 55559812 -> 52635740
$ws.Cells.Item(62, 9).Value = 71432140  # I62: 76926744 -> 71432140
$ws.Cells.Item(62, 11).Value = 71432140  # K62: 76926744 -> 71432140
$ws.Cells.Item(62, 13).Value = -71431516  # M62: -76926120 -> -71431516

# Row 65 (CRP)
$ws.Cells.Item(65, 8).Value = 52635740  # H65: 55559812 -> 52635740
$ws.Cells.Item(65, 9).Value = 71432140  # I65: 76926744 -> 71432140
$ws.Cells.Item(65, 11).Value = 357160700  # K65: 384633720 -> 357160700
$ws.Cells.Item(65, 13).Value = -357157580  # M65: -384630600 -> -357157580

# Row 96 (CRP)
$ws.Cells.Item(96, 8).Value = 15060.25  # H96: 13048.2 -> 15060.25
$ws.Cells.Item(96, 10).Value = 15060.25  # J96: 13048.2 -> 15060.25
$ws.Cells.Item(96, 12).Value = 15060.25  # L96: 13048.2 -> 15060.25
$ws.Cells.Item(96, 14).Value = -20552.25  # N96: -18540.2 -> -20552.25

# Row 134 (CRP)
$ws.Cells.Item(134, 8).Value = 999.2941  # H134: 953.1429000000001 -> 999.2941
$ws.Cells.Item(134, 9).Value = 849.7406999999999  # I134: 847.08 -> 849.7406999999999
$ws.Cells.Item(134, 10).Value = 1576.1428  # J134: 1218.3 -> 1576.1428
$ws.Cells.Item(134, 11).Value = 2549.2221  # K134: 2541.24 -> 2549.2221
$ws.Cells.Item(134, 12).Value = 4728.428400000001  # L134: 3654.9 -> 4728.428400000001
$ws.Cells.Item(134, 13).Value = -14.22209999999995  # M134: -6.240000000000236 -> -14.22209999999995
$ws.Cells.Item(134, 14).Value = -9798.428400000001  # N134: -8724.9 -> -9798.428400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Cells.Item(5, 8).Value = 1136.2632  # H5: 1238.625 -> 1136.2632
$ws.Cells.Item(5, 9).Value = 575.7  # I5: 659.625 -> 575.7
$ws.Cells.Item(5, 10).Value = 1759.1111  # J5: 1817.625 -> 1759.1111
$ws.Cells.Item(5, 11).Value = 1727.1  # K5: 1978.875 -> 1727.1
$ws.Cells.Item(5, 12).Value = 5277.3333  # L5: 5452.875 -> 5277.3333
$ws.Cells.Item(5, 13).Value = -1615.1  # M5: -1866.875 -> -1615.1
$ws.Cells.Item(5, 14).Value = -5501.3333  # N5: -5676.875 -> -5501.3333

# Row 131 (CUL)
$ws.Cells.Item(131, 8).Value = 753.6900000000001  # H131: 147842.77 -> 753.6900000000001
$ws.Cells.Item(131, 9).Value = 0  # I131: 1030 -> 0
$ws.Cells.Item(131, 10).Value = 753.6900000000001  # J131: 150034 -> 753.6900000000001
$ws.Cells.Item(131, 11).Value = 0  # K131: 3090 -> 0
$ws.Cells.Item(131, 12).Value = 2261.07  # L131: 450102 -> 2261.07
$ws.Cells.Item(131, 13).Value = $null  # M131: 1950 -> (removed)
$ws.Cells.Item(131, 14).Value = -12341.07  # N131: -460182 -> -12341.07

# Row 135 (CUL)
$ws.Cells.Item(135, 8).Value = 1136.2632  # H135: 1238.625 -> 1136.2632
$ws.Cells.Item(135, 9).Value = 575.7  # I135: 659.625 -> 575.7
$ws.Cells.Item(135, 10).Value = 1759.1111  # J135: 1817.625 -> 1759.1111
$ws.Cells.Item(135, 11).Value = 5181.3  # K135: 5936.625 -> 5181.3
$ws.Cells.Item(135, 12).Value = 15831.9999  # L135: 16358.625 -> 15831.9999
$ws.Cells.Item(135, 13).Value = -2646.3  # M135: -3401.625 -> -2646.3
$ws.Cells.Item(135, 14).Value = -20901.9999  # N135: -21428.625 -> -20901.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Cells.Item(80, 8).Value = 3250.2  # H80: 3789.95 -> 3250.2
$ws.Cells.Item(80, 9).Value = 2990  # I80: 3300 -> 2990
$ws.Cells.Item(80, 10).Value = 3423.6667  # J80: 4116.5835 -> 3423.6667
$ws.Cells.Item(80, 11).Value = 2990  # K80: 3300 -> 2990
$ws.Cells.Item(80, 12).Value = 3423.6667  # L80: 4116.5835 -> 3423.6667
$ws.Cells.Item(80, 13).Value = -1992  # M80: -2302 -> -1992
$ws.Cells.Item(80, 14).Value = -5419.6667  # N80: -6112.5835 -> -5419.6667

# Row 83 (GSM)
$ws.Cells.Item(83, 8).Value = 3250.2  # H83: 3789.95 -> 3250.2
$ws.Cells.Item(83, 9).Value = 2990  # I83: 3300 -> 2990
$ws.Cells.Item(83, 10).Value = 3423.6667  # J83: 4116.5835 -> 3423.6667
$ws.Cells.Item(83, 11).Value = 14950  # K83: 16500 -> 14950
$ws.Cells.Item(83, 12).Value = 17118.3335  # L83: 20582.9175 -> 17118.3335
$ws.Cells.Item(83, 13).Value = -9958  # M83: -11508 -> -9958
$ws.Cells.Item(83, 14).Value = -27102.3335  # N83: -30566.9175 -> -27102.3335

# Row 102 (GSM)
$ws.Cells.Item(102, 8).Value = 17858656  # H102: 15626381 -> 17858656
$ws.Cells.Item(102, 9).Value = 20834542  # I102: 17858244 -> 20834542
$ws.Cells.Item(102, 11).Value = 20834542  # K102: 17858244 -> 20834542
$ws.Cells.Item(102, 13).Value = -20832920  # M102: -17856622 -> -20832920

# Row 113 (GSM)
$ws.Cells.Item(113, 8).Value = 3510.282  # H113: 3593.5 -> 3510.282
$ws.Cells.Item(113, 9).Value = 4558.88  # I113: 4734.3335 -> 4558.88
$ws.Cells.Item(113, 11).Value = 4558.88  # K113: 4734.3335 -> 4558.88
$ws.Cells.Item(113, 13).Value = -2388.88  # M113: -2564.3335 -> -2388.88

# Row 123 (GSM)
$ws.Cells.Item(123, 8).Value = 5477.4546  # H123: 8656.625 -> 5477.4546
$ws.Cells.Item(123, 9).Value = 3121.4285  # I123: 4620 -> 3121.4285
$ws.Cells.Item(123, 10).Value = 9600.5  # J123: 10002.167 -> 9600.5
$ws.Cells.Item(123, 11).Value = 3121.4285  # K123: 4620 -> 3121.4285
$ws.Cells.Item(123, 12).Value = 9600.5  # L123: 10002.167 -> 9600.5
$ws.Cells.Item(123, 13).Value = -671.4285  # M123: -2170 -> -671.4285
$ws.Cells.Item(123, 14).Value = -14500.5  # N123: -14902.167 -> -14500.5

# Row 126 (GSM)
$ws.Cells.Item(126, 8).Value = 2577.15  # H126: 2947.9546 -> 2577.15
$ws.Cells.Item(126, 9).Value = 1582.0714  # I126: 1897.4 -> 1582.0714
$ws.Cells.Item(126, 10).Value = 4899  # J126: 5199.143 -> 4899
$ws.Cells.Item(126, 11).Value = 4746.2142  # K126: 5692.200000000001 -> 4746.2142
$ws.Cells.Item(126, 12).Value = 14697  # L126: 15597.429 -> 14697
$ws.Cells.Item(126, 13).Value = -2276.2142  # M126: -3222.200000000001 -> -2276.2142
$ws.Cells.Item(126, 14).Value = -19637  # N126: -20537.429 -> -19637

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 22487.926  # H132: 14980.683 -> 22487.926
$ws.Cells.Item(132, 9).Value = 4058  # I132: 2708.3784 -> 4058
$ws.Cells.Item(132, 10).Value = 103579.6  # J132: 128499.5 -> 103579.6
$ws.Cells.Item(132, 11).Value = 12174  # K132: 8125.135200000001 -> 12174
$ws.Cells.Item(132, 12).Value = 310738.8  # L132: 385498.5 -> 310738.8
$ws.Cells.Item(132, 13).Value = -9644  # M132: -5595.135200000001 -> -9644
$ws.Cells.Item(132, 14).Value = -315798.8  # N132: -390558.5 -> -315798.8

$ws = $wb.Worksheets.Item("LTW")
# Row 95 (LTW)
$ws.Cells.Item(95, 8).Value = 0  # H95: 40000 -> 0
$ws.Cells.Item(95, 10).Value = 0  # J95: 40000 -> 0
$ws.Cells.Item(95, 12).Value = 0  # L95: 40000 -> 0
$ws.Cells.Item(95, 14).Value = $null  # N95: -45492 -> (removed)

# Row 100 (LTW)
$ws.Cells.Item(100, 8).Value = 1840.2  # H100: 1764.8823 -> 1840.2
$ws.Cells.Item(100, 9).Value = 1278.4445  # I100: 1264.1818 -> 1278.4445
$ws.Cells.Item(100, 11).Value = 1278.4445  # K100: 1264.1818 -> 1278.4445
$ws.Cells.Item(100, 13).Value = -737.4445000000001  # M100: -723.1818000000001 -> -737.4445000000001

# Row 136 (LTW)
$ws.Cells.Item(136, 8).Value = 1407.5186  # H136: 1441 -> 1407.5186
$ws.Cells.Item(136, 9).Value = 1273.3636  # I136: 1314.4546 -> 1273.3636
$ws.Cells.Item(136, 11).Value = 3820.0908  # K136: 3943.3638 -> 3820.0908
$ws.Cells.Item(136, 13).Value = -1270.0908  # M136: -1393.3638 -> -1270.0908
